$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5
$ws.Range("A2").Value = 45052.50694444445
$ws.Range("B2").Value = 5.378
$ws.Range("C2").Value = 5.209
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 13.173
$ws.Range("F2").Value = 11.188
$ws.Range("G2").Value = 5.041
$ws.Range("H2").Value = 13.381
$ws.Range("I2").Value = 7.973
$ws.Range("J2").Value = 3.774
$ws.Range("K2").Value = 6.167
$ws.Range("L2").Value = 6.371
$ws.Range("M2").Value = 6.162
$ws.Range("N2").Value = 1.68
$ws.Range("O2").Value = 4.805
$ws.Range("P2").Value = 6.849
$ws.Range("Q2").Value = 4.447
$ws.Range("R2").Value = 0.297
$ws.Range("S2").Value = 0.425
$ws.Range("T2").Value = 72.19199999999999
$ws.Range("U2").Value = 14.308
$ws.Range("V2").Value = 4.642
$ws.Range("W2").Value = 8.505000000000001
$ws.Range("X2").Value = 6.25
$ws.Range("Y2").Value = 0.6899999999999999
$ws.Range("Z2").Value = 6.977
$ws.Range("AA2").Value = 4.068
$ws.Range("AB2").Value = 5.066
$ws.Range("AC2").Value = 6.321
$ws.Range("AD2").Value = 6.735
$ws.Range("AE2").Value = 0.773
$ws.Range("AF2").Value = 11.064
$ws.Range("AG2").Value = 3.76
$ws.Range("AH2").Value = 5.189
$ws.Range("A3").Value = 45052.51388888889
$ws.Range("B3").Value = 8.369999999999999
$ws.Range("C3").Value = 6.806
$ws.Range("D3").Value = 0.081
$ws.Range("E3").Value = 19.166
$ws.Range("F3").Value = 15.796
$ws.Range("G3").Value = 6.994
$ws.Range("H3").Value = 25.227
$ws.Range("I3").Value = 11.061
$ws.Range("J3").Value = 5.174
$ws.Range("K3").Value = 7.886
$ws.Range("L3").Value = 8.32
$ws.Range("M3").Value = 8.337999999999999
$ws.Range("N3").Value = 2.294
$ws.Range("O3").Value = 6.803
$ws.Range("P3").Value = 9.973000000000001
$ws.Range("Q3").Value = 5.958
$ws.Range("R3").Value = 0.186
$ws.Range("S3").Value = 0.414
$ws.Range("T3").Value = 101.101
$ws.Range("U3").Value = 19.925
$ws.Range("V3").Value = 6.43
$ws.Range("W3").Value = 12.892
$ws.Range("X3").Value = 7.69
$ws.Range("Y3").Value = 0.9360000000000001
$ws.Range("Z3").Value = 12.814
$ws.Range("AA3").Value = 5.75
$ws.Range("AB3").Value = 5.871
$ws.Range("AC3").Value = 6.955
$ws.Range("AD3").Value = 8.754
$ws.Range("AE3").Value = 0.419
$ws.Range("AF3").Value = 22.81
$ws.Range("AG3").Value = 4.173
$ws.Range("AH3").Value = 7.863
$ws.Range("A4").Value = 45052.52083333334
$ws.Range("B4").Value = 17.718
$ws.Range("C4").Value = 13.637
$ws.Range("D4").Value = 0.437
$ws.Range("E4").Value = 39.259
$ws.Range("F4").Value = 32.363
$ws.Range("G4").Value = 14.218
$ws.Range("H4").Value = 53.039
$ws.Range("I4").Value = 22.152
$ws.Range("J4").Value = 10.101
$ws.Range("K4").Value = 15.187
$ws.Range("L4").Value = 16.207
$ws.Range("M4").Value = 16.769
$ws.Range("N4").Value = 4.588
$ws.Range("O4").Value = 14.016
$ws.Range("P4").Value = 20.302
$ws.Range("Q4").Value = 11.856
$ws.Range("R4").Value = 0.23
$ws.Range("S4").Value = 0.641
$ws.Range("T4").Value = 209.93
$ws.Range("U4").Value = 39.896
$ws.Range("V4").Value = 13.056
$ws.Range("W4").Value = 26.619
$ws.Range("X4").Value = 14.609
$ws.Range("Y4").Value = 1.873
$ws.Range("Z4").Value = 26.368
$ws.Range("AA4").Value = 11.612
$ws.Range("AB4").Value = 10.761
$ws.Range("AC4").Value = 12.66
$ws.Range("AD4").Value = 17.042
$ws.Range("AE4").Value = 0.281
$ws.Range("AF4").Value = 48.109
$ws.Range("AG4").Value = 7.778
$ws.Range("AH4").Value = 16.245
$ws.Range("A5").Value = 45052.52777777778
$ws.Range("B5").Value = 24.08
$ws.Range("C5").Value = 18.33
$ws.Range("D5").Value = 0.68
$ws.Range("E5").Value = 52.96
$ws.Range("F5").Value = 43.69
$ws.Range("G5").Value = 19.16
$ws.Range("H5").Value = 74.08
$ws.Range("I5").Value = 29.72
$ws.Range("J5").Value = 13.48
$ws.Range("K5").Value = 20.19
$ws.Range("L5").Value = 21.61
$ws.Range("M5").Value = 22.54
$ws.Range("N5").Value = 6.16
$ws.Range("O5").Value = 18.95
$ws.Range("P5").Value = 27.35
$ws.Range("Q5").Value = 15.9
$ws.Range("R5").Value = 0.27
$ws.Range("S5").Value = 0.8
$ws.Range("T5").Value = 284.39
$ws.Range("U5").Value = 53.62
$ws.Range("V5").Value = 17.59
$ws.Range("W5").Value = 36.02
$ws.Range("X5").Value = 19.38
$ws.Range("Y5").Value = 2.51
$ws.Range("Z5").Value = 36.31
$ws.Range("AA5").Value = 15.61
$ws.Range("AB5").Value = 14.14
$ws.Range("AC5").Value = 16.62
$ws.Range("AD5").Value = 22.72
$ws.Range("AE5").Value = 0.21
$ws.Range("AF5").Value = 67.29000000000001
$ws.Range("AG5").Value = 10.28
$ws.Range("AH5").Value = 21.95

# Delete row 6 entirely
$ws.Rows.Item(6).Delete()

# Update column widths
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
